$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 735310.4
$ws.Range("J17").Value = 826974.5600000001
$ws.Range("L17").Value = 2480923.68
$ws.Range("N17").Value = -2481259.68
$ws.Range("H80").Value = 429.78946
$ws.Range("I80").Value = 418.42856
$ws.Range("J80").Value = 436.41666
$ws.Range("K80").Value = 1255.28568
$ws.Range("L80").Value = 1309.24998
$ws.Range("M80").Value = -257.28568
$ws.Range("N80").Value = -3305.24998
$ws.Range("H83").Value = 429.78946
$ws.Range("I83").Value = 418.42856
$ws.Range("J83").Value = 436.41666
$ws.Range("K83").Value = 3765.85704
$ws.Range("L83").Value = 3927.74994
$ws.Range("M83").Value = 1226.14296
$ws.Range("N83").Value = -13911.74994
$ws.Range("H116").Value = 4208.6924
$ws.Range("I116").Value = 4201.1816
$ws.Range("K116").Value = 4201.1816
$ws.Range("M116").Value = -759.1815999999999
$ws.Range("H133").Value = 100000
$ws.Range("J133").Value = 100000
$ws.Range("L133").Value = 100000
$ws.Range("N133").Value = -110120
$ws.Range("H135").Value = 88235670
$ws.Range("J135").Value = 500000300
$ws.Range("L135").Value = 4500002700
$ws.Range("N135").Value = -4500007770
$ws.Range("H141").Value = 5650
$ws.Range("I141").Value = 5650
$ws.Range("K141").Value = 16950
$ws.Range("M141").Value = -11770

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1698.25
$ws.Range("I45").Value = 1622
$ws.Range("K45").Value = 1622
$ws.Range("M45").Value = -1245
$ws.Range("H61").Value = 45456370
$ws.Range("I61").Value = 52633140
$ws.Range("K61").Value = 52633140
$ws.Range("M61").Value = -52632928
$ws.Range("H95").Value = 35344.11
$ws.Range("J95").Value = 35344.11
$ws.Range("L95").Value = 35344.11
$ws.Range("N95").Value = -40836.11
$ws.Range("H97").Value = 275.33334
$ws.Range("I97").Value = 282.70587
$ws.Range("J97").Value = 150
$ws.Range("K97").Value = 282.70587
$ws.Range("L97").Value = 150
$ws.Range("M97").Value = 213.29413
$ws.Range("N97").Value = -1142
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = None
$ws.Range("N105").ClearContents()
$ws.Range("H106").Value = 79999
$ws.Range("J106").Value = 79999
$ws.Range("L106").Value = 79999
$ws.Range("N106").Value = -82523
$ws.Range("H136").Value = 45456370
$ws.Range("I136").Value = 52633140
$ws.Range("K136").Value = 157899420
$ws.Range("M136").Value = -157896870
$ws.Range("H138").Value = 167498.75
$ws.Range("J138").Value = 167498.75
$ws.Range("L138").Value = 167498.75
$ws.Range("N138").Value = -177778.75
$ws.Range("H139").Value = 171959.2
$ws.Range("J139").Value = 171959.2
$ws.Range("L139").Value = 171959.2
$ws.Range("N139").Value = -182239.2

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2009.4
$ws.Range("I20").Value = 1869.4
$ws.Range("K20").Value = 1869.4
$ws.Range("M20").Value = -1622.4
$ws.Range("H107").Value = 57725.223
$ws.Range("I107").Value = 2256.8
$ws.Range("J107").Value = 127060.75
$ws.Range("K107").Value = 2256.8
$ws.Range("L107").Value = 127060.75
$ws.Range("M107").Value = -336.8000000000002
$ws.Range("N107").Value = -130900.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 5068.857
$ws.Range("I22").Value = 5068.857
$ws.Range("K22").Value = 5068.857
$ws.Range("M22").Value = -4718.857
$ws.Range("H25").Value = 1000
$ws.Range("I25").Value = 1000
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 1000
$ws.Range("L25").Value = None
$ws.Range("N25").ClearContents()
$ws.Range("M25").Value = -826
$ws.Range("H132").Value = 50001364
$ws.Range("I132").Value = 50001364
$ws.Range("K132").Value = 150004092
$ws.Range("M132").Value = -150001562
$ws.Range("H134").Value = 14708075
$ws.Range("I134").Value = 22729444
$ws.Range("J134").Value = 2229.6667
$ws.Range("K134").Value = 68188332
$ws.Range("L134").Value = 6689.000100000001
$ws.Range("M134").Value = -68185797
$ws.Range("N134").Value = -11759.0001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 34.42857
$ws.Range("I2").Value = 21.833334
$ws.Range("J2").Value = 43.875
$ws.Range("K2").Value = 131.000004
$ws.Range("L2").Value = 263.25
$ws.Range("M2").Value = -18.00000399999999
$ws.Range("N2").Value = -489.25
$ws.Range("H86").Value = 549.65216
$ws.Range("I86").Value = 403.53845
$ws.Range("J86").Value = 739.6
$ws.Range("K86").Value = 1210.61535
$ws.Range("L86").Value = 2218.8
$ws.Range("M86").Value = -24.61535000000003
$ws.Range("N86").Value = -4590.8
$ws.Range("H89").Value = 549.65216
$ws.Range("I89").Value = 403.53845
$ws.Range("J89").Value = 739.6
$ws.Range("K89").Value = 3631.84605
$ws.Range("L89").Value = 6656.400000000001
$ws.Range("M89").Value = 2296.15395
$ws.Range("N89").Value = -18512.4

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6911.2666
$ws.Range("J102").Value = 26275
$ws.Range("L102").Value = 26275
$ws.Range("N102").Value = -29519
$ws.Range("H122").Value = 7993.909
$ws.Range("I122").Value = 4462
$ws.Range("K122").Value = 13386
$ws.Range("M122").Value = -10936

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2388
$ws.Range("J7").Value = 1991
$ws.Range("L7").Value = 1991
$ws.Range("N7").Value = -2215
$ws.Range("H40").Value = 3050
$ws.Range("I40").Value = 3050
$ws.Range("K40").Value = 3050
$ws.Range("M40").Value = -2914
$ws.Range("H93").Value = 2048.8
$ws.Range("J93").Value = 4093.8
$ws.Range("L93").Value = 4093.8
$ws.Range("N93").Value = -6589.8
$ws.Range("H126").Value = 2388
$ws.Range("J126").Value = 1991
$ws.Range("L126").Value = 5973
$ws.Range("N126").Value = -10913
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = None
$ws.Range("N127").ClearContents()
$ws.Range("H132").Value = 20845802
$ws.Range("I132").Value = 25013964
$ws.Range("K132").Value = 75041892
$ws.Range("M132").Value = -75039362
$ws.Range("H136").Value = 1926.3043
$ws.Range("I136").Value = 1616.125
$ws.Range("K136").Value = 4848.375
$ws.Range("M136").Value = -2298.375

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 19742.5
$ws.Range("J45").Value = 19742.5
$ws.Range("L45").Value = 19742.5
$ws.Range("N45").Value = -20724.5
$ws.Range("H75").Value = 106999.75
$ws.Range("I75").Value = 28000
$ws.Range("K75").Value = 28000
$ws.Range("M75").Value = -27064
$ws.Range("H78").Value = 106999.75
$ws.Range("I78").Value = 28000
$ws.Range("K78").Value = 84000
$ws.Range("M78").Value = -79320
$ws.Range("H81").Value = 1299
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 1299
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H136").Value = 25002094
$ws.Range("I136").Value = 26317810
$ws.Range("K136").Value = 78953430
$ws.Range("M136").Value = -78950880

